{"js": "// The edit re-shuffles several blocks of text between fixed paragraph\n// \"slots\" in the document: the paragraph count/order/styles do NOT\n// change, only the w:t content inside certain paragraphs (and a few\n// runs inside the \"Avalia\u00e7\u00e3o\" paragraph) moves to a different slot.\n// We snapshot every source value first and only then write the new\n// values, so earlier writes can never corrupt a value still needed by\n// a later read.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst pObjPt = paragraphs.items[5];    // \"Apresentar aos estudantes...\" (under Objetivos)\nconst pObjEn = paragraphs.items[6];    // \"To introduce to students...\" (italic)\nconst pDocente = paragraphs.items[8];  // \"2342277 - Bertha Mar\u00eda Cuadros Melgar\" (ListBullet)\nconst pResumoPt = paragraphs.items[10]; // \"Carga e for\u00e7a el\u00e9trica, Campo el\u00e9trico...\"\nconst pResumoEn = paragraphs.items[11]; // \"Electric Charge and Matter...\" (italic)\nconst pProgramaPt = paragraphs.items[13]; // \"1) Carga e For\u00e7a el\u00e9trica: ...11) Equa\u00e7\u00f5es de Maxwell.\"\nconst pAvaliacao = paragraphs.items[16];  // M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o block\nconst pBibliografia = paragraphs.items[18]; // \"NUSSENZVEIG, H.M. ...JEWETT Jr...\"\n\npObjPt.load(\"text\");\npObjEn.load(\"text\");\npDocente.load(\"text\");\npResumoPt.load(\"text\");\npResumoEn.load(\"text\");\npProgramaPt.load(\"text\");\npBibliografia.load(\"text\");\nawait context.sync();\n\n// Snapshot every value currently in place before any mutation happens.\nconst objPt = pObjPt.text;\nconst objEn = pObjEn.text;\nconst docenteBullet = pDocente.text;\nconst resumoPt = pResumoPt.text;\nconst resumoEn = pResumoEn.text;\nconst programaPt = pProgramaPt.text;\nconst bibliografia = pBibliografia.text;\n\nconst metodoAnswer = \"NF=A avalia\u00e7\u00e3o ser\u00e1 composta por provas, listas, projetos, semin\u00e1rios e outras formas que far\u00e3o a composi\u00e7\u00e3o das notas, sendo estipulada a m\u00e9dia final a somat\u00f3ria destas notas (N), com no m\u00ednimo duas avalia\u00e7\u00f5es, sendo: (N1+...+Nn)/n.\";\nconst criterioAnswer = \"NF\u2265 5,0.\";\nconst normaAnswer = \"(NF+RC)/2 \u2265 5,0, onde RC \u00e9 uma prova de recupera\u00e7\u00e3o a ser aplicada.\";\n\n// --- Whole-paragraph slots -------------------------------------------------\npObjPt.insertText(resumoPt, Word.InsertLocation.replace);\npObjEn.insertText(resumoEn, Word.InsertLocation.replace);\npDocente.insertText(objPt, Word.InsertLocation.replace);\npResumoPt.insertText(programaPt, Word.InsertLocation.replace);\npResumoEn.insertText(objEn, Word.InsertLocation.replace);\npProgramaPt.insertText(metodoAnswer, Word.InsertLocation.replace);\npBibliografia.insertText(docenteBullet, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Sub-run answers inside the \"Avalia\u00e7\u00e3o\" paragraph ----------------------\n// The paragraph mixes bold \"Label:\" runs with plain-text answer runs; only\n// the three answer runs move, each one slot along the M\u00e9todo -> Crit\u00e9rio ->\n// Norma de recupera\u00e7\u00e3o -> Bibliografia chain. Replacing from the end of the\n// chain backwards means every search below still matches exactly one run\n// (the target hasn't been given the same text as an earlier link yet).\nasync function replaceAnswer(oldText, newText) {\n  const range = pAvaliacao.getRange();\n  const results = range.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for \" + JSON.stringify(oldText) +\n      \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceAnswer(normaAnswer, bibliografia);\nawait replaceAnswer(criterioAnswer, normaAnswer);\nawait replaceAnswer(metodoAnswer, criterioAnswer);\n", "ps1": "# The edit re-shuffles several blocks of text between fixed paragraph\n# \"slots\" in the document: the paragraph count/order/styles do NOT\n# change, only the text inside certain paragraphs (and a few runs\n# inside the \"Avalia\u00e7\u00e3o\" paragraph) moves to a different slot.\n# We snapshot every source value first and only then write the new\n# values, so earlier writes can never corrupt a value still needed by\n# a later read.\n\n$d = $word.ActiveDocument\n\n$pObjPt       = $d.Paragraphs.Item(6)    # \"Apresentar aos estudantes...\" (under Objetivos)\n$pObjEn       = $d.Paragraphs.Item(7)    # \"To introduce to students...\" (italic)\n$pDocente     = $d.Paragraphs.Item(9)    # \"2342277 - Bertha Mar\u00eda Cuadros Melgar\" (ListBullet)\n$pResumoPt    = $d.Paragraphs.Item(11)   # \"Carga e for\u00e7a el\u00e9trica, Campo el\u00e9trico...\"\n$pResumoEn    = $d.Paragraphs.Item(12)   # \"Electric Charge and Matter...\" (italic)\n$pProgramaPt  = $d.Paragraphs.Item(14)   # \"1) Carga e For\u00e7a el\u00e9trica: ...11) Equa\u00e7\u00f5es de Maxwell.\"\n$pAvaliacao   = $d.Paragraphs.Item(17)   # M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o block\n$pBibliografia = $d.Paragraphs.Item(19)  # \"NUSSENZVEIG, H.M. ...JEWETT Jr...\"\n\n# Snapshot every value currently in place before any mutation happens.\n# Paragraph.Range.Text includes the trailing paragraph-mark character\n# (CR, chr 13); strip it so re-assigning the captured text elsewhere\n# doesn't splice in an extra paragraph break.\n$objPt = $pObjPt.Range.Text.TrimEnd([char]13)\n$objEn = $pObjEn.Range.Text.TrimEnd([char]13)\n$docenteBullet = $pDocente.Range.Text.TrimEnd([char]13)\n$resumoPt = $pResumoPt.Range.Text.TrimEnd([char]13)\n$resumoEn = $pResumoEn.Range.Text.TrimEnd([char]13)\n$programaPt = $pProgramaPt.Range.Text.TrimEnd([char]13)\n$bibliografia = $pBibliografia.Range.Text.TrimEnd([char]13)\n\n$metodoAnswer = \"NF=A avalia\u00e7\u00e3o ser\u00e1 composta por provas, listas, projetos, semin\u00e1rios e outras formas que far\u00e3o a composi\u00e7\u00e3o das notas, sendo estipulada a m\u00e9dia final a somat\u00f3ria destas notas (N), com no m\u00ednimo duas avalia\u00e7\u00f5es, sendo: (N1+...+Nn)/n.\"\n$criterioAnswer = \"NF\u2265 5,0.\"\n$normaAnswer = \"(NF+RC)/2 \u2265 5,0, onde RC \u00e9 uma prova de recupera\u00e7\u00e3o a ser aplicada.\"\n\n# --- Whole-paragraph slots -------------------------------------------------\n# Paragraph.Range.Text includes the trailing paragraph mark in its \"length\"\n# but assigning to it only replaces the visible text, same as before.\n$pObjPt.Range.Text = $resumoPt\n$pObjEn.Range.Text = $resumoEn\n$pDocente.Range.Text = $objPt\n$pResumoPt.Range.Text = $programaPt\n$pResumoEn.Range.Text = $objEn\n$pProgramaPt.Range.Text = $metodoAnswer\n$pBibliografia.Range.Text = $docenteBullet\n\n# --- Sub-run answers inside the \"Avalia\u00e7\u00e3o\" paragraph ----------------------\n# The paragraph mixes bold \"Label:\" runs with plain-text answer runs; only\n# the three answer runs move, each one slot along the M\u00e9todo -> Crit\u00e9rio ->\n# Norma de recupera\u00e7\u00e3o -> Bibliografia chain. Replacing from the end of the\n# chain backwards means every Find below still matches exactly one run (the\n# target hasn't been given the same text as an earlier link yet). Find is\n# scoped to $pAvaliacao.Range so it can never match text in other\n# paragraphs.\n\nfunction Replace-AnswerInParagraph($para, [string]$oldText, [string]$newText) {\n    $r = $para.Range\n    $found = $r.Find.Execute($oldText, $true)\n    if (-not $found) {\n        throw \"text not found: $oldText\"\n    }\n    $r.Text = $newText\n}\n\nReplace-AnswerInParagraph $pAvaliacao $normaAnswer $bibliografia\nReplace-AnswerInParagraph $pAvaliacao $criterioAnswer $normaAnswer\nReplace-AnswerInParagraph $pAvaliacao $metodoAnswer $criterioAnswer\n"}
